# ExcelHelper: read Image column / ImageConverter: check null
#
# The commit adds an "Image" column (header + a file path value) to
# Sheet1, sizes the new column to fit its content on both sheets, and
# leaves the final on-screen selection on Sheet1!C4 (Sheet2's own
# selection becomes C1:C2, matching where its new column would line
# up once the Image column is added there too).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: add the new "Image" column ------------------------------
$ws1.Range("C1").Value = "Image"
$ws1.Range("C2").Value = "C:\Users\admin\Desktop\bmw.jpg"

# Size column C on both sheets to fit the new content.
$ws1.Columns.Item(3).ColumnWidth = 28.5
$ws2.Columns.Item(3).ColumnWidth = 28.5

# --- Selection bookkeeping --------------------------------------------
# Set Sheet2's own selection first (without leaving it the active tab).
$null = $ws2.Range("C1:C2").Select()

# Finish with Sheet1 active and C4 selected.
$ws1.Activate()
$null = $ws1.Range("C4").Select()
